$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row (row 11): Right and Wrong counts
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Update "Total" row (row 12): Right and Wrong counts, and the score fraction text
$ws.Range("B12").Value = 117
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "107/252"
